$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.243.35'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.807.72'
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '702.32'
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.53'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("D7").Value = '3.805.95'
$ws.Range("E7").Value = '  -1.01%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.49'
$ws.Range("E11").Value = '  +2.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.473'
$ws.Range("E12").Value = '  +3.34%  '
$ws.Range("E13").Value = '  -1.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.92'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").Value = '4.450.27'
$ws.Range("E15").Value = '  -1.02%  '
$ws.Range("D16").Value = '3.821.11'
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").Value = '71.296.69'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.115'
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.44'
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '514.61'
$ws.Range("E21").Value = '  +4.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.56'
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.714'
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.00'
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000142'
$ws.Range("E25").Value = '  -3.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.15'
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = '3.957.47'
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.33'
$ws.Range("E28").Value = '  -2.46%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.02'
$ws.Range("E30").Value = '  -3.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.02'
$ws.Range("E31").Value = '  -5.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.36'
$ws.Range("E32").Value = '  -1.86%  '
$ws.Range("E33").Value = '  -1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.05'
$ws.Range("E34").Value = '  -1.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.176'
$ws.Range("E35").Value = '  -2.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.12'
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("D37").Value = '3.770.24'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("E39").Value = '  -2.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.36'
$ws.Range("E40").Value = '  +5.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.37'
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.27'
$ws.Range("E43").Value = '  -1.79%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '173.86'
$ws.Range("E45").Value = '  +6.69%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000309'
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.38'
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '422.62'
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("E50").Value = '  -1.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.52'
$ws.Range("E51").Value = '  -0.99%  '
